# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (interest count) and "最低票价" (min ticket price)
# values to the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("G3").Value = 70
$ws1.Range("G4").Value = 80

$ws1.Range("F8").Value = 1767
$ws1.Range("F9").Value = 49
$ws1.Range("F11").Value = 148
$ws1.Range("F12").Value = 2021
$ws1.Range("F13").Value = 14
$ws1.Range("F14").Value = 139
$ws1.Range("F15").Value = 1028
$ws1.Range("F16").Value = 458

$ws1.Range("F25").Value = 6
$ws1.Range("F26").Value = 1089
$ws1.Range("F27").Value = 6
$ws1.Range("F28").Value = 328
$ws1.Range("F30").Value = 266
$ws1.Range("F31").Value = 314

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("G3").Value = 70
$ws4.Range("G4").Value = 80

$ws4.Range("F8").Value = 1767
$ws4.Range("F10").Value = 49
$ws4.Range("F12").Value = 148
$ws4.Range("F13").Value = 2021
$ws4.Range("F14").Value = 14
$ws4.Range("F15").Value = 139
$ws4.Range("F16").Value = 1028
$ws4.Range("F17").Value = 458

$ws4.Range("F26").Value = 6
$ws4.Range("F27").Value = 1089
$ws4.Range("F28").Value = 6
$ws4.Range("F29").Value = 328
$ws4.Range("F31").Value = 266
$ws4.Range("F32").Value = 314
